# Apply the OOXML diff to the report workbook.
#
# Sheet layout (see xl/workbook.xml): 1 = Overview, 2 = Issues, 3 = AI Summary
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Overview" ---------------------------------------------------
$overview = $wb.Worksheets.Item(1)

# B4: refreshed timestamp
$overview.Range("B4").Value = "2025-11-11T01:35:29.265000Z"

# B6: reworded AI summary
$overview.Range("B6").Value = "The site has a missing meta description, which can negatively impact its visibility in search results. Addressing this issue will improve click-through rates and overall SEO performance."

# --- Sheet 3: "AI Summary" --------------------------------------------------
$aiSummary = $wb.Worksheets.Item(3)

# Column widths: col A 56 -> 58, col C 64 -> 40 (col B stays at 100)
# The xlsx <col width="..."> attribute is offset from the COM ColumnWidth
# property by 5/6 (default Calibri 11 padding), so subtract that offset.
$aiSummary.Columns.Item(1).ColumnWidth = 58 - 5/6
$aiSummary.Columns.Item(3).ColumnWidth = 40 - 5/6

# B2: reworded summary (same new text as Overview!B6)
$aiSummary.Range("B2").Value = "The site has a missing meta description, which can negatively impact its visibility in search results. Addressing this issue will improve click-through rates and overall SEO performance."

# B8 / C8: reworded "why it matters" / evidence text
$aiSummary.Range("B8").Value = "Meta descriptions help search engines understand page content and influence user clicks."
$aiSummary.Range("C8").Value = "The homepage lacks a meta description."

# A14 / A15: reworded recommended actions
$aiSummary.Range("A14").Value = "Monitor click-through rates after implementing changes"
$aiSummary.Range("A15").Value = "Educate the team on best practices for meta descriptions"
